# Convert the opening "Heading1 title" + "bold byline" pair of paragraphs
# into a pandoc-style title block: a "Title"-styled paragraph (title split
# into word / space / word runs) followed by an "Authors"-styled paragraph
# (author name split the same way, with the leading "By " and bold removed).

$d = $word.ActiveDocument

$titlePara  = $d.Paragraphs(1)
$authorPara = $d.Paragraphs(2)

$rng = $d.Range($titlePara.Range.Start, $authorPara.Range.End)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml = "<w:p $wNs>" +
            "<w:pPr><w:pStyle w:val=`"Title`"/></w:pPr>" +
            "<w:r><w:t xml:space=`"preserve`">Beyond</w:t></w:r>" +
            "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
            "<w:r><w:t xml:space=`"preserve`">Politics</w:t></w:r>" +
          "</w:p>" +
          "<w:p $wNs>" +
            "<w:pPr><w:pStyle w:val=`"Authors`"/></w:pPr>" +
            "<w:r><w:t xml:space=`"preserve`">Dorothy</w:t></w:r>" +
            "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
            "<w:r><w:t xml:space=`"preserve`">Day</w:t></w:r>" +
          "</w:p>"

$rng.InsertXML($newXml)
